$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$centerAlign = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
$leftAlign   = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft
$pasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---------------------------------------------------------------------------
# Header row (row 1): rename "status" -> "status_penawaran" and
# "batas_atas_penawaran" -> "harga_penawaran"; center-align the three
# header cells that aren't the first column.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "status_penawaran"
$ws.Range("C1").Value = "harga_penawaran"
$ws.Range("B1:D1").HorizontalAlignment = $centerAlign

# ---------------------------------------------------------------------------
# "otomatis" block (rows 2-5) keeps its row position, but the wording of
# column B changes, and column B becomes left-aligned instead of centered.
# ---------------------------------------------------------------------------
$otomatis = @(
    @{ Row=2; B="kurang dari harga penawaran terakhir"; C=98000000; D="fail" },
    @{ Row=3; B="sama dengan harga penawaran terakhir"; C=98500000; D="fail" },
    @{ Row=4; B="lebih dari harga penawaran terakhir";  C=99000000; D="pass" },
    @{ Row=5; B="kurang dari harga minimum";            C=70000000; D="fail" }
)

foreach ($item in $otomatis) {
    $r = $item.Row
    $ws.Range("A" + $r).Value = "otomatis"
    $ws.Range("B" + $r).Value = $item.B
    $ws.Range("C" + $r).Value = $item.C
    $ws.Range("D" + $r).Value = $item.D
    $ws.Range("B" + $r).HorizontalAlignment = $leftAlign
}

# ---------------------------------------------------------------------------
# "manual" block (rows 6-9, new). Copy the formatting already established
# for the "otomatis" rows above (border + alignment) row-by-row, then
# overwrite values, then give columns A & D (jenis bid / expected) the
# distinct Arial/theme-colored font used for the manual block.
# ---------------------------------------------------------------------------
$manual = @(
    @{ Row=6; SrcRow=2; B="kurang dari harga penawaran terakhir"; C=100500000;  D="fail" },
    @{ Row=7; SrcRow=3; B="sama dengan harga penawaran terakhir"; C=2000000000; D="fail" },
    @{ Row=8; SrcRow=4; B="lebih dari harga penawaran terakhir";  C=2005000000; D="pass" },
    @{ Row=9; SrcRow=5; B="kurang dari harga minimum";            C=50000000;   D="fail" }
)

foreach ($item in $manual) {
    $r = $item.Row
    $src = $item.SrcRow

    $ws.Range("A" + $src + ":D" + $src).Copy()
    $ws.Range("A" + $r + ":D" + $r).PasteSpecial($pasteFormats)

    $ws.Range("A" + $r).Value = "manual"
    $ws.Range("B" + $r).Value = $item.B
    $ws.Range("C" + $r).Value = $item.C
    $ws.Range("D" + $r).Value = $item.D

    $ws.Range("A" + $r).Font.Name = "Arial"
    $ws.Range("A" + $r).Font.ThemeColor = 1
    $ws.Range("D" + $r).Font.Name = "Arial"
    $ws.Range("D" + $r).Font.ThemeColor = 1
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Column B needs to be wider to fit the longer status text (target raw
# OOXML width ~= 34.14 characters; 33.25 is the closest COM ColumnWidth
# input that round-trips to that value given Excel's internal
# character/pixel snapping).
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 33.25
